# Updated symbol list refresh: new Price (column D) and Volume(1h) (column E)
# quotes for the coin rows. Values are entered with a leading apostrophe so
# that Excel stores them as literal text (matching the workbook's existing
# inline-string cells) instead of auto-converting numeric-looking strings
# (e.g. "245.88") or percentages (e.g. "-0.10%") into Number-typed cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.88"
$ws.Range("E2").Value = "'-0.10%"

$ws.Range("D3").Value = "'28.27"
$ws.Range("E3").Value = "'-3.59%"

$ws.Range("D4").Value = "'5.290"
$ws.Range("E4").Value = "'1.94%"

$ws.Range("D5").Value = "'0.05713"

$ws.Range("D6").Value = "'6.647"
$ws.Range("E6").Value = "'1.19%"

$ws.Range("E7").Value = "'3.51%"

$ws.Range("D8").Value = "'0.8637"
$ws.Range("E8").Value = "'0.52%"

$ws.Range("D9").Value = "'0.8847"
$ws.Range("E9").Value = "'2.47%"

$ws.Range("D10").Value = "'0.1388"
$ws.Range("E10").Value = "'1.67%"

$ws.Range("D11").Value = "'0.07081"
$ws.Range("E11").Value = "'-0.11%"

$ws.Range("D12").Value = "'0.03118"
$ws.Range("E12").Value = "'1.94%"

$ws.Range("D13").Value = "'0.09225"
$ws.Range("E13").Value = "'-1.61%"

$ws.Range("D14").Value = "'0.001529"
$ws.Range("E14").Value = "'-0.41%"

$ws.Range("D15").Value = "'0.0005975"
$ws.Range("E15").Value = "'-0.28%"

$ws.Range("D16").Value = "'0.005981"
$ws.Range("E16").Value = "'0.29%"

$ws.Range("D17").Value = "'3.494"
$ws.Range("E17").Value = "'-0.02%"

$ws.Range("E18").Value = "'-0.65%"

$ws.Range("D19").Value = "'0.3168"
$ws.Range("E19").Value = "'-0.86%"

$ws.Range("D20").Value = "'0.03343"
$ws.Range("E20").Value = "'1.29%"

$ws.Range("E21").Value = "'0.68%"

$ws.Range("D22").Value = "'3.489"
$ws.Range("E22").Value = "'0.50%"

$ws.Range("D23").Value = "'0.04102"
$ws.Range("E23").Value = "'-1.34%"

$ws.Range("E24").Value = "'-0.19%"

$ws.Range("D25").Value = "'0.001223"
$ws.Range("E25").Value = "'-0.36%"

$ws.Range("D26").Value = "'0.004169"
$ws.Range("E26").Value = "'-16.49%"

$ws.Range("E40").Value = "'1.04%"

$ws.Range("E41").Value = "'-0.46%"

$ws.Range("D42").Value = "'0.002198"
$ws.Range("E42").Value = "'-8.82%"

$ws.Range("D43").Value = "'0.002948"
$ws.Range("E43").Value = "'-16.28%"

$ws.Range("D44").Value = "'0.009464"
$ws.Range("E44").Value = "'11.82%"

$ws.Range("D45").Value = "'0.00005274"
$ws.Range("E45").Value = "'-0.16%"

$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.11%"

$ws.Range("D47").Value = "'0.08903"
$ws.Range("E47").Value = "'56.13%"

$ws.Range("E48").Value = "'-0.34%"

$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.11%"

$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.11%"
